$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Service ID values (kept as text to match original shared-string typing)
$idCells = @("C2","C3","C4","C5","C11","C12","C13","C14","C24")
foreach ($addr in $idCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C2").Value = "10234686"
$ws.Range("C3").Value = "10234687"
$ws.Range("C4").Value = "10234689"
$ws.Range("C5").Value = "10234691"
$ws.Range("C11").Value = "10234708"
$ws.Range("C12").Value = "10234710"
$ws.Range("C13").Value = "10234717"
$ws.Range("C14").Value = "10234718"
$ws.Range("C24").Value = "136546157"

# Updated Selenium failure log text (shared by F12, F13, F14)
$errText = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=120.0.6099.110)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.37', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '19.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 120.0.6099.110, chrome: {chromedriverVersion: 120.0.6099.109 (3419140ab66..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:62982}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 38ad931a7ef8dd60159bc4e7d91db684
*** Element info: {Using=id, value=lblServiceID}
'@

$ws.Range("F12").Value = $errText
$ws.Range("F13").Value = $errText
$ws.Range("F14").Value = $errText

# Restore default row height after auto-fit triggered by multi-line text
$ws.Rows(12).RowHeight = 15
$ws.Rows(13).RowHeight = 15
$ws.Rows(14).RowHeight = 15
